$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:G3").NumberFormat = "@"

$ws.Range("A3").Value = "gf"
$ws.Range("B3").Value = "44"
$ws.Range("C3").Value = "fdg"
$ws.Range("D3").Value = "df"
$ws.Range("E3").Value = "5634653546546"
$ws.Range("F3").Value = "555"
$ws.Range("G3").Value = "fgg"

$ws.Range("A3:G3").Style = "Normal"
